# Edit Assignment3/ProjectManagment/actions/actions7.docx:
# 1) Nudge/resize the horizontal rule shape (Shape1) and flip it vertically
#    (positionV, extent, xfrm ext + flipV, and the VML fallback from/to + flip:y).
# 2) Add the text "Brandon allocated to " to the previously-empty run in the
#    last (numbered) paragraph.

$d = $word.ActiveDocument

# --- Part 1: update the line-shape drawing (paragraph 2 holds Shape1) ---
# The shape's own geometry (flipV, independent xfrm extent, VML fallback
# coordinates) isn't reachable through the high-level Shape.* properties in
# this host, so we replace the run's content with the target OOXML via
# Range.InsertXML, keeping every other part of the run untouched.
$shapeParaXml = '<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" mc:Ignorable="w14 wp14"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Noto Sans Cond" w:hAnsi="Noto Sans Cond"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Noto Sans Cond" w:hAnsi="Noto Sans Cond"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor behindDoc="0" distT="17780" distB="17780" distL="17780" distR="17780" simplePos="0" locked="0" layoutInCell="0" allowOverlap="1" relativeHeight="2" wp14:anchorId="6039D491"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>259715</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>262890</wp:posOffset></wp:positionV><wp:extent cx="5804535" cy="3810"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="1" name="Shape1"/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm flipV="1"><a:off x="0" y="0"/><a:ext cx="5803920" cy="1800"/></a:xfrm><a:prstGeom prst="line"><a:avLst/></a:prstGeom><a:ln w="36360"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:round/></a:ln></wps:spPr><wps:style><a:lnRef idx="0"></a:lnRef><a:fillRef idx="0"/><a:effectRef idx="0"></a:effectRef><a:fontRef idx="minor"/></wps:style><wps:bodyPr/></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:line id="shape_0" from="20.45pt,20.7pt" to="477.4pt,20.8pt" ID="Shape1" stroked="t" style="position:absolute;flip:y" wp14:anchorId="6039D491"><v:stroke color="black" weight="36360" joinstyle="round" endcap="flat"/><v:fill o:detectmouseclick="t" on="false"/><w10:wrap type="none"/></v:line></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$shapePara = $d.Paragraphs.Item(2)
$shapePara.Range.InsertXML($shapeParaXml)

# --- Part 2: add text to the empty run in the last paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
# Exclude the trailing paragraph mark so the text lands inside the existing run.
$insertRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$insertRange.InsertAfter("Brandon allocated to ")

Write-Output "edit applied"
